# hiyerarşi ve yorum düzenlemeleri - TABLO 2 icin Quiz sutunu kaldirildi,
# TOPLAM (G) sutunu silindi ve degerler guncellendi.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quiz sutunu (D) kaldirildiginda toplamlar yeniden hesaplaniyor; once eski
# TOPLAM sutununu (G) sil, boylece G1:G8 verileri ve F-oncesi bicimler kayar.
$ws.Columns("G:G").Delete()

# Baslik satiri - F1 artik yeni TOPLAM basligi
$ws.Range("F1").Value = "TOPLAM"

# Satir 2 (ust ozet satiri) guncel degerleri
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 40
$ws.Range("F2").ClearContents()

# Satir 3 - basliklar kaydi (Quiz kaldirildi)
$ws.Range("D3").Value = "Vize"
$ws.Range("E3").Value = "Fin"
$ws.Range("F3").Value = "TOPLAM"

# Satir 4
$ws.Range("D4").Value = 0
$ws.Range("F4").Value = 1

# Satir 5
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 3

# Satir 6
$ws.Range("D6").Value = 1
$ws.Range("F6").Value = 3

# Satir 7
$ws.Range("D7").Value = 1
$ws.Range("F7").Value = 2

# Satir 8
$ws.Range("F8").Value = 2
